$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 18:38"

# --- Re-rank three countries: Pakistan moves up above Chile/Irlanda ---
# Row 26 becomes Pakistan (updated stats), row 27 becomes Chile (former row26
# data, country unchanged), row 28 becomes Irlanda (former row27 data,
# country unchanged).
$ws.Range("A26").Value = "Pakistan"
$ws.Range("B26").Value = 22048
$ws.Range("C26").Value = 1107
$ws.Range("D26").Value = 5801
$ws.Range("E26").Value = 15733
$ws.Range("F26").Value = 111
$ws.Range("G26").Value = 38
$ws.Range("H26").Value = 514

$ws.Range("A27").Value = "Chile"
$ws.Range("B27").Value = 22016
$ws.Range("C27").Value = 1373
$ws.Range("D27").Value = 10710
$ws.Range("E27").Value = 11031
$ws.Range("F27").Value = 470
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 275

$ws.Range("A28").Value = "Irlanda"
$ws.Range("B28").Value = 21772
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 13386
$ws.Range("E28").Value = 7067
$ws.Range("F28").Value = 93
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 1319

# --- Plain statistic refreshes (country/rank unchanged) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1218638
$ws.Range("C4").Value = 5803
$ws.Range("D4").Value = 188778
$ws.Range("E4").Value = 959133
$ws.Range("F4").Value = 16055
$ws.Range("G4").Value = 806
$ws.Range("H4").Value = 70727

# Row 7: Reino Unido
$ws.Range("B7").Value = 194990
$ws.Range("C7").Value = 4406
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 165219
$ws.Range("F7").Value = 1559
$ws.Range("G7").Value = 693
$ws.Range("H7").Value = 29427

# Row 11: Turquia
$ws.Range("B11").Value = 129491
$ws.Range("C11").Value = 1832
$ws.Range("D11").Value = 73285
$ws.Range("E11").Value = 52686
$ws.Range("F11").Value = 1338
$ws.Range("G11").Value = 59
$ws.Range("H11").Value = 3520

# Row 48: Chequia
$ws.Range("B48").Value = 7878
$ws.Range("C48").Value = 59
$ws.Range("D48").Value = 3999
$ws.Range("E48").Value = 3625
$ws.Range("F48").Value = 59
$ws.Range("G48").Value = 2
$ws.Range("H48").Value = 254

# Row 56: Marruecos
$ws.Range("B56").Value = 5219
$ws.Range("C56").Value = 166
$ws.Range("D56").Value = 1838
$ws.Range("E56").Value = 3200
$ws.Range("F56").Value = 1
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 181

# Row 77: Islandia
$ws.Range("D77").Value = 1733
$ws.Range("E77").Value = 56

# Row 143: Madagascar
$ws.Range("B143").Value = 151
$ws.Range("C143").Value = 2
$ws.Range("D143").Value = 101
